# Update '想去人数' (want-to-go count) values per the diff, across all 4 sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 42293
$ws.Range("F3").Value = 27
$ws.Range("F4").Value = 10013
$ws.Range("F6").Value = 1044
$ws.Range("F7").Value = 972
$ws.Range("F8").Value = 777
$ws.Range("F9").Value = 242
$ws.Range("F10").Value = 316
$ws.Range("F11").Value = 1014
$ws.Range("F14").Value = 800
$ws.Range("F15").Value = 348
$ws.Range("F16").Value = 1615
$ws.Range("F18").Value = 795
$ws.Range("F21").Value = 720
$ws.Range("F22").Value = 800
$ws.Range("F26").Value = 573
$ws.Range("F27").Value = 572
$ws.Range("F28").Value = 77
$ws.Range("F31").Value = 25
$ws.Range("F32").Value = 457
$ws.Range("F35").Value = 175
$ws.Range("F37").Value = 1425
$ws.Range("F39").Value = 1315
$ws.Range("F41").Value = 108
$ws.Range("F45").Value = 53
$ws.Range("F46").Value = 15

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 350
$ws.Range("F8").Value = 157
$ws.Range("F16").Value = 41

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 2103
$ws.Range("F3").Value = 569
$ws.Range("F4").Value = 473

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 2103
$ws.Range("F3").Value = 569
$ws.Range("F4").Value = 350
$ws.Range("F5").Value = 27
$ws.Range("F6").Value = 10013
$ws.Range("F7").Value = 1044
$ws.Range("F8").Value = 1044
$ws.Range("F10").Value = 473
$ws.Range("F11").Value = 973
$ws.Range("F12").Value = 777
$ws.Range("F13").Value = 157
$ws.Range("F14").Value = 316
$ws.Range("F15").Value = 1014
$ws.Range("F17").Value = 800
$ws.Range("F18").Value = 348
$ws.Range("F19").Value = 1615
$ws.Range("F21").Value = 795
$ws.Range("F24").Value = 720
$ws.Range("F25").Value = 800
$ws.Range("F29").Value = 573
$ws.Range("F31").Value = 572
$ws.Range("F32").Value = 77
$ws.Range("F36").Value = 25
$ws.Range("F37").Value = 457
$ws.Range("F40").Value = 175
$ws.Range("F43").Value = 1315
